$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "PartOfSponza" (1st sheet) - new column F (v1243) step-perf
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PartOfSponza")
$ws1.Activate()

$ws1.Cells.Item(1, 6).Value = "v1243"

$ws1F = @(102, 102, 103, 101, 103, 101, 102, 102, 102, 101)
for ($i = 0; $i -lt $ws1F.Length; $i++) {
    $ws1.Cells.Item($i + 2, 6).Value = $ws1F[$i]
}

$ws1.Range("F12").Select()

# ------------------------------------------------------------------
# Sheet "Sponza" (2nd sheet) - new column C (v1243) total-perf
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sponza")
$ws2.Activate()

$ws2.Cells.Item(1, 3).Value = "v1243"

$ws2C = @(8844, 8917, 8703, 9020, 9197, 8943, 9048, 9011, 9055, 8898)
for ($i = 0; $i -lt $ws2C.Length; $i++) {
    $ws2.Cells.Item($i + 2, 3).Value = $ws2C[$i]
}

$ws2.Range("C14").Select()

# ------------------------------------------------------------------
# Sheet "ComplexMesh" (3rd sheet) - new column C (v1243) total-perf
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ComplexMesh")
$ws3.Activate()

$ws3.Cells.Item(1, 3).Value = "v1243"

$ws3C = @(5553, 5561, 5600, 5607, 5577, 5598, 5617, 5560, 5563, 5573)
for ($i = 0; $i -lt $ws3C.Length; $i++) {
    $ws3.Cells.Item($i + 2, 3).Value = $ws3C[$i]
}

$ws3.Range("D9").Select()

# ComplexMesh ends up the active / selected tab, matching the saved view.
$ws3.Activate()
